# Fix code indentation in module (slide 39, "Group 7" -> "TextBox 4")
#
# The Julia "module" code sample had its inner statements (export, the two
# function headers, and their matching "end"s) at the same indentation as
# the enclosing "module"/"end" keywords. This indents those five lines by
# two spaces, leaving the outer "module ..." / final "end" untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(39)
$grp = $s.Shapes.Item(4)          # "Group 7"
$tb = $grp.GroupItems.Item(1)     # "TextBox 4" (the Julia module snippet)
$tr = $tb.TextFrame.TextRange

# Paragraph 2: "export fact"          -> "  export fact"
$tr.Paragraphs(2, 1).Runs(1, 1).Text = "  export "

# Paragraph 4: "function fact(n)"     -> "  function fact(n)"
$tr.Paragraphs(4, 1).Runs(1, 1).Text = "  function "

# Paragraph 6: "end" (closes fact)    -> "  end"
$tr.Paragraphs(6, 1).Runs(1, 1).Text = "  end"

# Paragraph 8: "function gcd(a, b)"   -> "  function gcd(a, b)"
$tr.Paragraphs(8, 1).Runs(1, 1).Text = "  function "

# Paragraph 10: "end" (closes gcd)    -> "  end"
$tr.Paragraphs(10, 1).Runs(1, 1).Text = "  end"

# Paragraph 12 ("end" closing the module) is left unchanged.
